$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 'Volume 32   Number  12'
$ws.Range("C9").Value = 'Report Covering the Week  3/17/2025  Through  3/23/2025'
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = -100
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = -66.666666666666
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 11
$ws.Range("K14").Value = -54.545454545454
$ws.Range("L14").Value = -64.285714285714
$ws.Range("M14").Value = -54.545454545454
$ws.Range("N14").Value = -93.333333333333
$ws.Range("C15").Value = 5
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 400
$ws.Range("F15").Value = 15
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 49
$ws.Range("J15").Value = 27
$ws.Range("K15").Value = 81.481481481481
$ws.Range("L15").Value = 58.064516129032
$ws.Range("M15").Value = 25.641025641025
$ws.Range("N15").Value = -51
$ws.Range("C16").Value = 25
$ws.Range("E16").Value = -26.470588235294
$ws.Range("F16").Value = 119
$ws.Range("G16").Value = 154
$ws.Range("H16").Value = -22.727272727272
$ws.Range("I16").Value = 366
$ws.Range("J16").Value = 450
$ws.Range("K16").Value = -18.666666666666
$ws.Range("L16").Value = -8.270676691729
$ws.Range("M16").Value = -26.060606060606
$ws.Range("N16").Value = -82.041216879293
$ws.Range("C17").Value = 50
$ws.Range("D17").Value = 67
$ws.Range("E17").Value = -25.373134328358
$ws.Range("F17").Value = 202
$ws.Range("G17").Value = 263
$ws.Range("H17").Value = -23.193916349809
$ws.Range("I17").Value = 561
$ws.Range("J17").Value = 659
$ws.Range("K17").Value = -14.871016691957
$ws.Range("L17").Value = -7.730263157894
$ws.Range("M17").Value = 44.587628865979
$ws.Range("N17").Value = -47.323943661971
$ws.Range("D18").Value = 27
$ws.Range("E18").Value = -3.703703703703
$ws.Range("F18").Value = 96
$ws.Range("G18").Value = 126
$ws.Range("H18").Value = -23.809523809523
$ws.Range("I18").Value = 289
$ws.Range("J18").Value = 331
$ws.Range("K18").Value = -12.688821752265
$ws.Range("L18").Value = -23.947368421052
$ws.Range("M18").Value = -11.076923076923
$ws.Range("N18").Value = -88.082474226804
$ws.Range("C19").Value = 138
$ws.Range("D19").Value = 142
$ws.Range("E19").Value = -2.81690140845
$ws.Range("F19").Value = 489
$ws.Range("G19").Value = 522
$ws.Range("H19").Value = -6.321839080459
$ws.Range("I19").Value = 1212
$ws.Range("J19").Value = 1397
$ws.Range("K19").Value = -13.242662848962
$ws.Range("L19").Value = -5.973622963537
$ws.Range("M19").Value = 24.819773429454
$ws.Range("N19").Value = -48.730964467005
$ws.Range("C20").Value = 24
$ws.Range("D20").Value = 18
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 60
$ws.Range("G20").Value = 70
$ws.Range("H20").Value = -14.285714285714
$ws.Range("I20").Value = 151
$ws.Range("J20").Value = 202
$ws.Range("K20").Value = -25.247524752475
$ws.Range("L20").Value = -40.551181102362
$ws.Range("M20").Value = 69.662921348314
$ws.Range("N20").Value = -92.93401965372
$ws.Range("C21").Value = 268
$ws.Range("D21").Value = 290
$ws.Range("E21").Value = -7.586206896551
$ws.Range("F21").Value = 982
$ws.Range("G21").Value = 1143
$ws.Range("H21").Value = -14.085739282589
$ws.Range("I21").Value = 2633
$ws.Range("J21").Value = 3077
$ws.Range("K21").Value = -14.429639259018
$ws.Range("L21").Value = -11.495798319327
$ws.Range("M21").Value = 13.589301121656
$ws.Range("N21").Value = -74.196393571148
$ws.Range("C22").Value = 3
$ws.Range("E22").Value = -40
$ws.Range("F22").Value = 14
$ws.Range("G22").Value = 22
$ws.Range("H22").Value = -36.363636363636
$ws.Range("I22").Value = 51
$ws.Range("J22").Value = 68
$ws.Range("K22").Value = -25
$ws.Range("L22").Value = -19.047619047619
$ws.Range("M22").Value = -5.555555555555
$ws.Range("C23").Value = 25
$ws.Range("D23").Value = 18
$ws.Range("E23").Value = 38.888888888888
$ws.Range("F23").Value = 91
$ws.Range("G23").Value = 101
$ws.Range("H23").Value = -9.900990099009
$ws.Range("I23").Value = 272
$ws.Range("J23").Value = 289
$ws.Range("K23").Value = -5.882352941176
$ws.Range("L23").Value = 1.492537313432
$ws.Range("M23").Value = 52.808988764044
$ws.Range("C24").Value = 272
$ws.Range("D24").Value = 245
$ws.Range("E24").Value = 11.020408163265
$ws.Range("F24").Value = 1062
$ws.Range("G24").Value = 949
$ws.Range("H24").Value = 11.90727081138
$ws.Range("I24").Value = 3180
$ws.Range("J24").Value = 2704
$ws.Range("K24").Value = 17.603550295858
$ws.Range("L24").Value = 4.743083003952
$ws.Range("M24").Value = 77.752934600335
$ws.Range("C25").Value = 164
$ws.Range("D25").Value = 123
$ws.Range("E25").Value = 33.333333333333
$ws.Range("G25").Value = 481
$ws.Range("H25").Value = 27.858627858627
$ws.Range("I25").Value = 1821
$ws.Range("J25").Value = 1419
$ws.Range("K25").Value = 28.329809725158
$ws.Range("L25").Value = 5.565217391304
$ws.Range("C26").Value = 83
$ws.Range("D26").Value = 91
$ws.Range("E26").Value = -8.791208791208
$ws.Range("F26").Value = 358
$ws.Range("G26").Value = 380
$ws.Range("H26").Value = -5.78947368421
$ws.Range("I26").Value = 968
$ws.Range("J26").Value = 1017
$ws.Range("K26").Value = -4.818092428711
$ws.Range("L26").Value = 1.361256544502
$ws.Range("M26").Value = -12.159709618874
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 500
$ws.Range("F27").Value = 19
$ws.Range("G27").Value = 11
$ws.Range("H27").Value = 72.727272727272
$ws.Range("I27").Value = 58
$ws.Range("J27").Value = 44
$ws.Range("K27").Value = 31.818181818181
$ws.Range("L27").Value = 0
$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 17
$ws.Range("E28").Value = -23.529411764705
$ws.Range("F28").Value = 48
$ws.Range("G28").Value = 55
$ws.Range("H28").Value = -12.727272727272
$ws.Range("I28").Value = 124
$ws.Range("J28").Value = 127
$ws.Range("K28").Value = -2.362204724409
$ws.Range("L28").Value = -0.8
$ws.Range("D29").Value = 1
$ws.Range("E29").Value = 100
$ws.Range("F29").Value = 6
$ws.Range("G29").Value = 7
$ws.Range("H29").Value = -14.285714285714
$ws.Range("I29").Value = 14
$ws.Range("J29").Value = 21
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = -58.823529411764
$ws.Range("M29").Value = -60
$ws.Range("N29").Value = -92.134831460674
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 100
$ws.Range("F30").Value = 6
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 13
$ws.Range("J30").Value = 17
$ws.Range("K30").Value = -23.529411764705
$ws.Range("L30").Value = -59.375
$ws.Range("M30").Value = -60.60606060606
$ws.Range("N30").Value = -92.121212121212
$ws.Range("D31").Value = 5
$ws.Range("G31").Value = 12
$ws.Range("H31").Value = -75
$ws.Range("J31").Value = 23
$ws.Range("K31").Value = -56.521739130434
$ws.Range("L31").Value = -44.444444444444
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("G33").Value = 2
$ws.Range("J33").Value = 2
$ws.Range("K33").Value = 100
$ws.Range("L33").Value = -33.333333333333
